# Fill in "Membre du groupe" (group member) roster onto Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row, header row, then the 4 group members (Prénom(s) / Nom).
$data = @(
    @("Membres du groupe",  ""),
    @("Prénom (s)",         "Nom"),
    @("Alioune Abdou Salam","Kane"),
    @("Awa",                "Diaw"),
    @("Ange Emilson Rayan", "Raherinasolo"),
    @("Khadidiatou",        "Diakhaté")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Merge the title across both columns and center the title + header rows.
$ws.Range("A1:B1").Merge()
$ws.Range("A2:B2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:B1").HorizontalAlignment = -4108   # xlCenter
